$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the payment code referenced in the test data (02-03-2023 -> 03-03-2023 run) ---
# Column F (PREPARATION) contains the full login/preparation script including the code.
$ws.Range("F2").Value = "Username : 37841;`nPassword : bni1234;`nRole : RL09 - Penyelia Settlement;`nKode Pembayaran : DISK230300012"
# Column N (KODE_PEMBAYARAN) holds just the payment code value.
$ws.Range("N2").Value = "DISK230300012"

# --- Realign a handful of data cells (B2, D2, E2, F2) from center to left horizontal alignment ---
foreach ($addr in @("B2", "D2", "E2", "F2")) {
    $ws.Range($addr).HorizontalAlignment = -4131
}

# --- Move/update the active selection on the sheet ---
$ws.Range("P2").Select()
